$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 53
$ws.Range("A53").Value = "12/4/2012"
$ws.Range("B53").Value = 2.5
$ws.Range("D53").Value = "Manual: new section continued, new figure for illustration"

# Row 54
$ws.Range("A54").Value = "12/5/2012"
$ws.Range("B54").Value = 0.75
$ws.Range("C54").Value = 2.5
$ws.Range("D54").Value = "Manual: new section 2.7 completed, including new figures and first review"

# Row 55
$ws.Range("A55").Value = "12/6/2012"
$ws.Range("B55").Value = 2.5
$ws.Range("D55").Value = "Manual: Section 4.4, data type system time rewritten"

# Update selection to match the new active cell
$ws.Range("D55").Select()
